$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 24: 1/31/2020 ---
$ws.Range("A24").NumberFormat = "@"
$ws.Range("A24").Value = "1/31/2020"
$ws.Range("B24").Value = 0.5375
$ws.Range("C24").Value = "None"
$ws.Range("D24").Value = "Review the course, and prepare for the assignments"
$ws.Range("E24").Value = "Finished Goal"
$ws.Range("F24").Value = "The features our team find last time are definitely not essential features. I thought that the Maintable in our project is an essential feature, but how to define the boundary of this feature might be a problem. "
$ws.Range("G24").Value = "Not bad"
$ws.Range("A23:G23").Copy()
$ws.Range("A24:G24").PasteSpecial(-4122)
$ws.Rows.Item(24).RowHeight = 115

# --- Row 25: 2/2/2020 ---
$ws.Range("A25").NumberFormat = "@"
$ws.Range("A25").Value = "2/2/2020"
$ws.Range("B25").Value = 0.891666666666667
$ws.Range("C25").Value = "None"
$ws.Range("D25").Value = "Discuss the assignments"
$ws.Range("E25").Value = "Finished Goal"
$ws.Range("F25").Value = "Our team decides to have a team discussion next week. I want to prepare something for the discussion, so I try to find some clues about the feature ""Maintable"""
$ws.Range("G25").Value = "Not bad"
$ws.Range("A23:G23").Copy()
$ws.Range("A25:G25").PasteSpecial(-4122)
$ws.Rows.Item(25).RowHeight = 100

# --- Row 26: 2/4/2020 ---
$ws.Range("A26").NumberFormat = "@"
$ws.Range("A26").Value = "2/4/2020"
$ws.Range("B26").Value = 0.954166666666667
$ws.Range("C26").Value = "None"
$ws.Range("D26").Value = "Team discussion"
$ws.Range("E26").Value = "Finished Goal"
$ws.Range("F26").Value = "We had the team discussion this afternoon, and divided tasks to different members. Although we had some divergence, we have agreement in the end. We set the MainTable and EntryEditor as our two features. Besides, we want to explain the MainTable in a kind of high level and EntryEditor in details. "
$ws.Range("G26").Value = "Not bad"
$ws.Range("A23:G23").Copy()
$ws.Range("A26:G26").PasteSpecial(-4122)
$ws.Rows.Item(26).RowHeight = 164

# --- Row 27: 2/5/2020 ---
$ws.Range("A27").NumberFormat = "@"
$ws.Range("A27").Value = "2/5/2020"
$ws.Range("B27").Value = 0.625
$ws.Range("C27").Value = "None"
$ws.Range("D27").Value = "Do assignments"
$ws.Range("E27").Value = "Almost done"
$ws.Range("F27").Value = "We discussed a lot of things about the assignment, such as the format, what kind of diagram we need, etc. We add more things to our assignment, it's almost done."
$ws.Range("G27").Value = "Not bad"
$ws.Range("A23:G23").Copy()
$ws.Range("A27:G27").PasteSpecial(-4122)
$ws.Rows.Item(27).RowHeight = 99

# --- Sheet view: scroll / selection position ---
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C27").Select()
